$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("xref_waypoints_utms")

# --- Update/replace rows 42-47 (previously just site/name stubs for
#     McDowell Creek ef waypoints, now split into numbered ef1..ef4 /
#     ef1..ef2 points each with wp/gps_correction_s/easting/northing) ---

$ws.Cells.Item(42,1).Value = "McDowell Creek"
$ws.Cells.Item(42,2).Value = "58159_ds_ef1"
$ws.Cells.Item(42,3).Value = 68
$ws.Cells.Item(42,4).Value = 1
$ws.Cells.Item(42,5).Value = 627533
$ws.Cells.Item(42,6).Value = 6060447

$ws.Cells.Item(43,1).Value = "McDowell Creek"
$ws.Cells.Item(43,2).Value = "58159_ds_ef2"
$ws.Cells.Item(43,3).Value = 69
$ws.Cells.Item(43,4).Value = 1
$ws.Cells.Item(43,5).Value = 627555
$ws.Cells.Item(43,6).Value = 6060461

$ws.Cells.Item(44,1).Value = "McDowell Creek"
$ws.Cells.Item(44,2).Value = "58159_ds_ef3"
$ws.Cells.Item(44,3).Value = 71
$ws.Cells.Item(44,4).Value = 1
$ws.Cells.Item(44,5).Value = 627585
$ws.Cells.Item(44,6).Value = 6060475

$ws.Cells.Item(45,1).Value = "McDowell Creek"
$ws.Cells.Item(45,2).Value = "58159_ds_ef4"
$ws.Cells.Item(45,3).Value = 72
$ws.Cells.Item(45,4).Value = 1
$ws.Cells.Item(45,5).Value = 627611
$ws.Cells.Item(45,6).Value = 6060455

$ws.Cells.Item(46,1).Value = "McDowell Creek"
$ws.Cells.Item(46,2).Value = "58159_us_ef1"
$ws.Cells.Item(46,3).Value = 74
$ws.Cells.Item(46,4).Value = 1
$ws.Cells.Item(46,5).Value = 627669
$ws.Cells.Item(46,6).Value = 6060433

$ws.Cells.Item(47,1).Value = "McDowell Creek"
$ws.Cells.Item(47,2).Value = "58159_us_ef2"
$ws.Cells.Item(47,3).Value = 75
$ws.Cells.Item(47,4).Value = 1
$ws.Cells.Item(47,5).Value = 627702
$ws.Cells.Item(47,6).Value = 6060433

# --- New rows 48-57 appended below the previous last row (47) ---

$ws.Cells.Item(48,1).Value = "McDowell Creek"
$ws.Cells.Item(48,2).Value = "58159_us_ef3"
$ws.Cells.Item(48,5).Value = 627793
$ws.Cells.Item(48,6).Value = 6060395

$ws.Cells.Item(49,1).Value = "McDowell Creek"
$ws.Cells.Item(49,2).Value = "58159_us_ef4"
$ws.Cells.Item(49,5).Value = 628003
$ws.Cells.Item(49,6).Value = 6060508

$ws.Cells.Item(50,1).Value = "Gibson Creek"
$ws.Cells.Item(50,2).Value = "195290_us_ef1"
$ws.Cells.Item(50,5).Value = 640028
$ws.Cells.Item(50,6).Value = 6051717

$ws.Cells.Item(51,1).Value = "Gibson Creek"
$ws.Cells.Item(51,2).Value = "195290_ds_ef1"
$ws.Cells.Item(51,5).Value = 639984
$ws.Cells.Item(51,6).Value = 6051683

$ws.Cells.Item(52,1).Value = "Gibson Creek"
$ws.Cells.Item(52,2).Value = "195290_us_ef2"
$ws.Cells.Item(52,3).Value = 97
$ws.Cells.Item(52,4).Value = 1
$ws.Cells.Item(52,5).Value = 640890
$ws.Cells.Item(52,6).Value = 6051588

$ws.Cells.Item(53,1).Value = "Johnny David Creek"
$ws.Cells.Item(53,2).Value = "197663_ds_ef1"
$ws.Cells.Item(53,3).Value = 66
$ws.Cells.Item(53,4).Value = 1
$ws.Cells.Item(53,5).Value = 670225
$ws.Cells.Item(53,6).Value = 6044638

$ws.Cells.Item(54,1).Value = "Johnny David Creek"
$ws.Cells.Item(54,2).Value = "197663_us_ef1"
$ws.Cells.Item(54,3).Value = 53
$ws.Cells.Item(54,4).Value = 1
$ws.Cells.Item(54,5).Value = 670225
$ws.Cells.Item(54,6).Value = 6044812

$ws.Cells.Item(55,1).Value = "Riddeck Creek"
$ws.Cells.Item(55,2).Value = "197360_us_ef1"
$ws.Cells.Item(55,3).Value = 268
$ws.Cells.Item(55,4).Value = 2
$ws.Cells.Item(55,5).Value = 650101
$ws.Cells.Item(55,6).Value = 5992660

$ws.Cells.Item(56,1).Value = "Riddeck Creek"
$ws.Cells.Item(56,2).Value = "197360_us_ef2"
$ws.Cells.Item(56,5).Value = 650179
$ws.Cells.Item(56,6).Value = 5992702

$ws.Cells.Item(57,1).Value = "Riddeck Creek"
$ws.Cells.Item(57,2).Value = "197360_ds_ef1"
$ws.Cells.Item(57,5).Value = 649896
$ws.Cells.Item(57,6).Value = 5992406

# --- Cosmetic: widen column B to fit the new longer "*_efN" labels,
#     matching the bestFit width added for column B in the edit ---
$ws.Columns.Item(2).ColumnWidth = 14

# --- View state: scroll down to the newly-added rows and leave the
#     selection where the author's cursor ended up ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("C62").Select()
